$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2022.7646
$ws.Range("I51").Value = 1547.5
$ws.Range("J51").Value = 2169
$ws.Range("K51").Value = 1547.5
$ws.Range("L51").Value = 2169
$ws.Range("M51").Value = -1063.5
$ws.Range("N51").Value = -3137
$ws.Range("H98").Value = 2576.25
$ws.Range("I98").Value = 2810.2903
$ws.Range("J98").Value = 1125.2
$ws.Range("K98").Value = 2810.2903
$ws.Range("L98").Value = 1125.2
$ws.Range("M98").Value = -1312.2903
$ws.Range("N98").Value = -4121.2
$ws.Range("H112").Value = 1823.875
$ws.Range("I112").Value = 811.25
$ws.Range("J112").Value = 2026.4
$ws.Range("K112").Value = 2433.75
$ws.Range("L112").Value = 6079.200000000001
$ws.Range("M112").Value = -1325.75
$ws.Range("N112").Value = -8295.200000000001
$ws.Range("H113").Value = 3490.6
$ws.Range("I113").Value = 3478.75
$ws.Range("K113").Value = 3478.75
$ws.Range("M113").Value = -224.75
$ws.Range("H116").Value = 2083.7827
$ws.Range("I116").Value = 1655.4615
$ws.Range("J116").Value = 2640.6
$ws.Range("K116").Value = 1655.4615
$ws.Range("L116").Value = 2640.6
$ws.Range("M116").Value = 1786.5385
$ws.Range("N116").Value = -9524.6
$ws.Range("H122").Value = 2576.25
$ws.Range("I122").Value = 2810.2903
$ws.Range("J122").Value = 1125.2
$ws.Range("K122").Value = 8430.8709
$ws.Range("L122").Value = 3375.6
$ws.Range("M122").Value = -5980.8709
$ws.Range("N122").Value = -8275.6
$ws.Range("H137").Value = 1278.4423
$ws.Range("I137").Value = 886.8125
$ws.Range("J137").Value = 1905.05
$ws.Range("K137").Value = 2660.4375
$ws.Range("L137").Value = 5715.15
$ws.Range("M137").Value = -110.4375
$ws.Range("N137").Value = -10815.15
$ws.Range("H138").Value = 1273.9196
$ws.Range("I138").Value = 806.4737
$ws.Range("J138").Value = 1636.4286
$ws.Range("K138").Value = 2419.4211
$ws.Range("L138").Value = 4909.2858
$ws.Range("M138").Value = 2720.5789
$ws.Range("N138").Value = -15189.2858
$ws.Range("H140").Value = 34386
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 34386
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 34386
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -44746

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3898.9136
$ws.Range("I32").Value = 3414.6711
$ws.Range("K32").Value = 3414.6711
$ws.Range("M32").Value = -3127.6711
$ws.Range("H61").Value = 55557344
$ws.Range("I61").Value = 76924664
$ws.Range("J61").Value = 2322.8
$ws.Range("K61").Value = 76924664
$ws.Range("L61").Value = 2322.8
$ws.Range("M61").Value = -76924452
$ws.Range("N61").Value = -2746.8
$ws.Range("H102").Value = 7937823
$ws.Range("I102").Value = 8773069
$ws.Range("J102").Value = 2981
$ws.Range("K102").Value = 8773069
$ws.Range("L102").Value = 2981
$ws.Range("M102").Value = -8771447
$ws.Range("N102").Value = -6225
$ws.Range("H122").Value = 1156.8889
$ws.Range("I122").Value = 1177.5217
$ws.Range("J122").Value = 1038.25
$ws.Range("K122").Value = 3532.5651
$ws.Range("L122").Value = 3114.75
$ws.Range("M122").Value = -1082.5651
$ws.Range("N122").Value = -8014.75
$ws.Range("H132").Value = 1740.381
$ws.Range("I132").Value = 1834.4667
$ws.Range("J132").Value = 1505.1666
$ws.Range("K132").Value = 5503.4001
$ws.Range("L132").Value = 4515.4998
$ws.Range("M132").Value = -2973.4001
$ws.Range("N132").Value = -9575.4998
$ws.Range("H136").Value = 55557344
$ws.Range("I136").Value = 76924664
$ws.Range("J136").Value = 2322.8
$ws.Range("K136").Value = 230773992
$ws.Range("L136").Value = 6968.400000000001
$ws.Range("M136").Value = -230771442
$ws.Range("N136").Value = -12068.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 66668736
$ws.Range("I105").Value = 100001820
$ws.Range("J105").Value = 2562.2
$ws.Range("K105").Value = 100001820
$ws.Range("L105").Value = 2562.2
$ws.Range("M105").Value = -100000073
$ws.Range("N105").Value = -6056.2
$ws.Range("H134").Value = 3251.1914
$ws.Range("I134").Value = 935.1795
$ws.Range("J134").Value = 14541.75
$ws.Range("K134").Value = 2805.5385
$ws.Range("L134").Value = 43625.25
$ws.Range("M134").Value = -270.5384999999997
$ws.Range("N134").Value = -48695.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 50001250
$ws.Range("I16").Value = 66667840
$ws.Range("J16").Value = 1465
$ws.Range("K16").Value = 66667840
$ws.Range("L16").Value = 1465
$ws.Range("M16").Value = -66667553
$ws.Range("N16").Value = -2039
$ws.Range("H31").Value = 1881.8529
$ws.Range("I31").Value = 1758.069
$ws.Range("J31").Value = 2599.8
$ws.Range("K31").Value = 1758.069
$ws.Range("L31").Value = 2599.8
$ws.Range("M31").Value = -1463.069
$ws.Range("N31").Value = -3189.8
$ws.Range("H34").Value = 1881.8529
$ws.Range("I34").Value = 1758.069
$ws.Range("J34").Value = 2599.8
$ws.Range("K34").Value = 1758.069
$ws.Range("L34").Value = 2599.8
$ws.Range("M34").Value = -1556.069
$ws.Range("N34").Value = -3003.8
$ws.Range("H86").Value = 1862848.1
$ws.Range("J86").Value = 20136.25
$ws.Range("L86").Value = 20136.25
$ws.Range("N86").Value = -22382.25
$ws.Range("H89").Value = 1862848.1
$ws.Range("J89").Value = 20136.25
$ws.Range("L89").Value = 100681.25
$ws.Range("N89").Value = -111913.25
$ws.Range("H107").Value = 557.4783
$ws.Range("I107").Value = 460.2857
$ws.Range("K107").Value = 460.2857
$ws.Range("M107").Value = 1459.7143
$ws.Range("H113").Value = 50001250
$ws.Range("I113").Value = 66667840
$ws.Range("J113").Value = 1465
$ws.Range("K113").Value = 66667840
$ws.Range("L113").Value = 1465
$ws.Range("M113").Value = -66665670
$ws.Range("N113").Value = -5805
$ws.Range("H132").Value = 3060.4167
$ws.Range("I132").Value = 3260.3208
$ws.Range("J132").Value = 1546.8572
$ws.Range("K132").Value = 9780.9624
$ws.Range("L132").Value = 4640.571599999999
$ws.Range("M132").Value = -7250.9624
$ws.Range("N132").Value = -9700.571599999999
$ws.Range("H134").Value = 10205296
$ws.Range("I134").Value = 1224.289
$ws.Range("J134").Value = 125001100
$ws.Range("K134").Value = 3672.867
$ws.Range("L134").Value = 375003300
$ws.Range("M134").Value = -1137.867
$ws.Range("N134").Value = -375008370

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 367.6154
$ws.Range("I5").Value = 339.91666
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 1019.74998
$ws.Range("L5").Value = 2100
$ws.Range("M5").Value = -907.7499799999999
$ws.Range("N5").Value = -2324
$ws.Range("H107").Value = 5358.95
$ws.Range("J107").Value = 8060.5386
$ws.Range("L107").Value = 24181.6158
$ws.Range("N107").Value = -28021.6158
$ws.Range("H122").Value = 755.0625
$ws.Range("I122").Value = 662.9231
$ws.Range("J122").Value = 818.1053000000001
$ws.Range("K122").Value = 5966.3079
$ws.Range("L122").Value = 7362.947700000001
$ws.Range("M122").Value = -3516.3079
$ws.Range("N122").Value = -12262.9477
$ws.Range("H131").Value = 76925544
$ws.Range("J131").Value = 5506.6
$ws.Range("L131").Value = 16519.8
$ws.Range("N131").Value = -26599.8
$ws.Range("H135").Value = 367.6154
$ws.Range("I135").Value = 339.91666
$ws.Range("J135").Value = 700
$ws.Range("K135").Value = 3059.24994
$ws.Range("L135").Value = 6300
$ws.Range("M135").Value = -524.2499399999997
$ws.Range("N135").Value = -11370
$ws.Range("H140").Value = 24905.307
$ws.Range("I140").Value = 49644.523
$ws.Range("J140").Value = 3020.6155
$ws.Range("K140").Value = 148933.569
$ws.Range("L140").Value = 9061.8465
$ws.Range("M140").Value = -143753.569
$ws.Range("N140").Value = -19421.8465

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 31253402
$ws.Range("I102").Value = 41668332
$ws.Range("K102").Value = 41668332
$ws.Range("M102").Value = -41666710
$ws.Range("H113").Value = 1694.5
$ws.Range("I113").Value = 1577.8572
$ws.Range("J113").Value = 1966.6666
$ws.Range("K113").Value = 1577.8572
$ws.Range("L113").Value = 1966.6666
$ws.Range("M113").Value = 592.1428000000001
$ws.Range("N113").Value = -6306.6666
$ws.Range("H122").Value = 1339.7858
$ws.Range("I122").Value = 1135.1538
$ws.Range("K122").Value = 3405.4614
$ws.Range("M122").Value = -955.4614000000001
$ws.Range("H126").Value = 2991.5386
$ws.Range("I126").Value = 1875
$ws.Range("K126").Value = 5625
$ws.Range("M126").Value = -3155
$ws.Range("H132").Value = 1689.8182
$ws.Range("I132").Value = 1371.1724
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 4113.5172
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -1583.5172
$ws.Range("N132").Value = -17060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1201.3529
$ws.Range("I16").Value = 1418.4166
$ws.Range("K16").Value = 1418.4166
$ws.Range("M16").Value = -1248.4166
$ws.Range("H61").Value = 1746
$ws.Range("I61").Value = 1416.5
$ws.Range("K61").Value = 1416.5
$ws.Range("M61").Value = -1214.5
$ws.Range("H113").Value = 1746
$ws.Range("I113").Value = 1416.5
$ws.Range("K113").Value = 1416.5
$ws.Range("M113").Value = 753.5
$ws.Range("H132").Value = 26815.125
$ws.Range("I132").Value = 1622.2963
$ws.Range("K132").Value = 4866.8889
$ws.Range("M132").Value = -2336.8889

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 573
$ws.Range("I107").Value = 504.85715
$ws.Range("K107").Value = 1514.57145
$ws.Range("M107").Value = 405.4285500000001
$ws.Range("H122").Value = 11908960
$ws.Range("I122").Value = 13162366
$ws.Range("J122").Value = 1604.5
$ws.Range("K122").Value = 39487098
$ws.Range("L122").Value = 4813.5
$ws.Range("M122").Value = -39484648
$ws.Range("N122").Value = -9713.5
$ws.Range("H132").Value = 2325.484
$ws.Range("I132").Value = 2997.8948
$ws.Range("K132").Value = 8993.6844
$ws.Range("M132").Value = -6463.6844
